$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: the Price/Volume/Hora columns hold numeric-looking text, so values
# are entered with a leading apostrophe to force Excel to store them as text
# (matching the workbook's existing inline-string convention) rather than
# auto-converting to numbers/percentages.

# Rows 2-5: price/volume/hora updates only (no name/link changes)
$ws.Range("D2").Value = "'296.67"
$ws.Range("E2").Value = "'3.15%"
$ws.Range("G2").Value = "'18"

$ws.Range("D3").Value = "'41.57"
$ws.Range("E3").Value = "'3.27%"
$ws.Range("G3").Value = "'18"

$ws.Range("D4").Value = "'5.053"
$ws.Range("E4").Value = "'0.30%"
$ws.Range("G4").Value = "'18"

$ws.Range("D5").Value = "'0.07516"
$ws.Range("E5").Value = "'3.32%"
$ws.Range("G5").Value = "'18"

# Rows 6-17: coin list shifted up by one rank, with updated price/volume figures.
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.580"
$ws.Range("E6").Value = "'4.56%"
$ws.Range("G6").Value = "'18"

$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9258"
$ws.Range("E7").Value = "'1.30%"
$ws.Range("G7").Value = "'18"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.425"
$ws.Range("E8").Value = "'1.17%"
$ws.Range("G8").Value = "'18"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1214"
$ws.Range("E9").Value = "'1.59%"
$ws.Range("G9").Value = "'18"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1835"
$ws.Range("E10").Value = "'6.94%"
$ws.Range("G10").Value = "'18"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08831"
$ws.Range("E11").Value = "'3.69%"
$ws.Range("G11").Value = "'18"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04143"
$ws.Range("E12").Value = "'-0.30%"
$ws.Range("G12").Value = "'18"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.12%"
$ws.Range("G13").Value = "'18"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001292"
$ws.Range("E14").Value = "'1.61%"
$ws.Range("G14").Value = "'18"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005851"
$ws.Range("E15").Value = "'-1.35%"
$ws.Range("G15").Value = "'18"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.344"
$ws.Range("E16").Value = "'-1.55%"
$ws.Range("G16").Value = "'18"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.357"
$ws.Range("E17").Value = "'1.77%"
$ws.Range("G17").Value = "'18"

# Rows 18-26: price/volume updates only
$ws.Range("D18").Value = "'0.3315"
$ws.Range("E18").Value = "'1.67%"
$ws.Range("G18").Value = "'18"

$ws.Range("D19").Value = "'7.893"
$ws.Range("E19").Value = "'1.39%"
$ws.Range("G19").Value = "'18"

$ws.Range("D20").Value = "'0.1408"
$ws.Range("E20").Value = "'4.07%"
$ws.Range("G20").Value = "'18"

$ws.Range("D21").Value = "'0.2962"
$ws.Range("E21").Value = "'2.67%"
$ws.Range("G21").Value = "'18"

$ws.Range("D22").Value = "'0.04040"
$ws.Range("E22").Value = "'5.02%"
$ws.Range("G22").Value = "'18"

$ws.Range("D23").Value = "'0.001262"
$ws.Range("G23").Value = "'18"

$ws.Range("D24").Value = "'0.003875"
$ws.Range("E24").Value = "'1.59%"
$ws.Range("G24").Value = "'18"

$ws.Range("D25").Value = "'0.0001228"
$ws.Range("E25").Value = "'-4.17%"
$ws.Range("G25").Value = "'18"

$ws.Range("E26").Value = "'-0.07%"
$ws.Range("G26").Value = "'18"

# Rows 27-37: only the Hora column changes (no price data for these coins)
$ws.Range("G27").Value = "'18"
$ws.Range("G28").Value = "'18"
$ws.Range("G29").Value = "'18"
$ws.Range("G30").Value = "'18"
$ws.Range("G31").Value = "'18"
$ws.Range("G32").Value = "'18"
$ws.Range("G33").Value = "'18"
$ws.Range("G34").Value = "'18"
$ws.Range("G35").Value = "'18"
$ws.Range("G36").Value = "'18"
$ws.Range("G37").Value = "'18"

# Rows 38-51: price/volume updates
$ws.Range("D38").Value = "'0.02409"
$ws.Range("E38").Value = "'5.01%"
$ws.Range("G38").Value = "'18"

$ws.Range("D39").Value = "'0.05207"
$ws.Range("E39").Value = "'5.38%"
$ws.Range("G39").Value = "'18"

$ws.Range("D40").Value = "'0.006206"
$ws.Range("E40").Value = "'-9.97%"
$ws.Range("G40").Value = "'18"

$ws.Range("D41").Value = "'0.007782"
$ws.Range("E41").Value = "'1.02%"
$ws.Range("G41").Value = "'18"

$ws.Range("D42").Value = "'0.1324"
$ws.Range("E42").Value = "'4.35%"
$ws.Range("G42").Value = "'18"

$ws.Range("D43").Value = "'0.007357"
$ws.Range("E43").Value = "'-0.21%"
$ws.Range("G43").Value = "'18"

$ws.Range("D44").Value = "'0.008091"
$ws.Range("E44").Value = "'6.29%"
$ws.Range("G44").Value = "'18"

$ws.Range("D45").Value = "'0.2976"
$ws.Range("E45").Value = "'-4.39%"
$ws.Range("G45").Value = "'18"

$ws.Range("D46").Value = "'0.00006243"
$ws.Range("E46").Value = "'-2.14%"
$ws.Range("G46").Value = "'18"

$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.28%"
$ws.Range("G47").Value = "'18"

$ws.Range("D48").Value = "'0.04643"
$ws.Range("E48").Value = "'-81.56%"
$ws.Range("G48").Value = "'18"

$ws.Range("D49").Value = "'0.004194"
$ws.Range("E49").Value = "'-0.19%"
$ws.Range("G49").Value = "'18"

$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.28%"
$ws.Range("G50").Value = "'18"

$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.28%"
$ws.Range("G51").Value = "'18"
